$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The observation previously on row 20 moves to row 21, and the
# observation previously on row 21 moves to row 20 (only the cells that
# actually differ between the two rows need to be touched).

# --- Save old values from row 20 ---
$A20 = $ws.Range("A20").Value2
$B20 = $ws.Range("B20").Value2
$E20 = $ws.Range("E20").Value2
$F20 = $ws.Range("F20").Value2
$G20 = $ws.Range("G20").Value2
$H20 = $ws.Range("H20").Value2
$P20 = $ws.Range("P20").Value2
$Q20 = $ws.Range("Q20").Value2
$R20 = $ws.Range("R20").Value2
$S20 = $ws.Range("S20").Value2
$Z20 = $ws.Range("Z20").Value2
$AB20 = $ws.Range("AB20").Value2
$AW20 = $ws.Range("AW20").Value2
$AX20 = $ws.Range("AX20").Value2

# --- Save old values from row 21 ---
$A21 = $ws.Range("A21").Value2
$B21 = $ws.Range("B21").Value2
$E21 = $ws.Range("E21").Value2
$F21 = $ws.Range("F21").Value2
$G21 = $ws.Range("G21").Value2
$H21 = $ws.Range("H21").Value2
$M21 = $ws.Range("M21").Value2
$P21 = $ws.Range("P21").Value2
$Q21 = $ws.Range("Q21").Value2
$R21 = $ws.Range("R21").Value2
$S21 = $ws.Range("S21").Value2
$Z21 = $ws.Range("Z21").Value2
$AB21 = $ws.Range("AB21").Value2
$AC21 = $ws.Range("AC21").Value2
$AW21 = $ws.Range("AW21").Value2
$AX21 = $ws.Range("AX21").Value2

# --- Write row 21's old values into row 20 ---
$ws.Range("A20").Value2 = $A21
$ws.Range("B20").Value2 = $B21
$ws.Range("E20").Value2 = $E21
$ws.Range("F20").Value2 = $F21
$ws.Range("G20").Value2 = $G21
$ws.Range("H20").Value2 = $H21
$ws.Range("M20").Value2 = $M21
$ws.Range("P20").Value2 = $P21
$ws.Range("Q20").Value2 = $Q21
$ws.Range("R20").Value2 = $R21
$ws.Range("S20").Value2 = $S21
$ws.Range("Z20").Value2 = $Z21
$ws.Range("AB20").Value2 = $AB21
$ws.Range("AC20").Value2 = $AC21
$ws.Range("AW20").Value2 = $AW21
$ws.Range("AX20").Value2 = $AX21

# --- Write row 20's old values into row 21 ---
$ws.Range("A21").Value2 = $A20
$ws.Range("B21").Value2 = $B20
$ws.Range("E21").Value2 = $E20
$ws.Range("F21").Value2 = $F20
$ws.Range("G21").Value2 = $G20
$ws.Range("H21").Value2 = $H20
$ws.Range("M21").Value2 = ""
$ws.Range("P21").Value2 = $P20
$ws.Range("Q21").Value2 = $Q20
$ws.Range("R21").Value2 = $R20
$ws.Range("S21").Value2 = $S20
$ws.Range("Z21").Value2 = $Z20
$ws.Range("AB21").Value2 = $AB20
$ws.Range("AC21").Value2 = ""
$ws.Range("AW21").Value2 = $AW20
$ws.Range("AX21").Value2 = $AX20
